$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($c = 26; $c -ge 21; $c--) {
    $ws.Cells.Item(3, $c).Copy($ws.Cells.Item(3, $c + 1))
    $ws.Cells.Item(4, $c).Copy($ws.Cells.Item(4, $c + 1))
}
$ws.Cells.Item(3, 20).Copy($ws.Cells.Item(3, 21))
$ws.Cells.Item(4, 20).Copy($ws.Cells.Item(4, 21))
$ws.Cells.Item(3, 21).Value2 = "Организация"
$ws.Cells.Item(4, 21).Value2 = '${e.organization}'

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:AA4"))

# Fix up the ListColumns metadata: position 27 is a placeholder "ColumnNN"
# that actually belongs, physically, at the very end (AA) where "Дата
# увольнения" now lives; shift all the in-between names down by one so the
# metadata list matches the real physical header row again.
for ($i = 27; $i -ge 22; $i--) {
    $lo.ListColumns.Item($i).Name = $lo.ListColumns.Item($i - 1).Name
}
$lo.ListColumns.Item(21).Name = "Организация"

for ($i=1; $i -le $lo.ListColumns.Count; $i++) {
    Write-Output ($i.ToString() + ": " + $lo.ListColumns.Item($i).Name)
}
